$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. In the original document it
#    sits right after the "Файербол ПКМ" run, near the top of the file.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Append a single space to the very end of the document (end of the
#    last paragraph, right before the final paragraph mark), keeping the
#    sz/szCs = 28 (14pt) formatting used throughout that paragraph, and
#    landing it in its own run (matching two separate <w:r> elements).
#
#    Using Find/Replace on the already-formatted trailing text makes the
#    replacement inherit the full run properties (including w:szCs,
#    which isn't reachable through Font.Size alone). Nudging the size of
#    just the trailing character away and back then forces the engine to
#    split that trailing character into its own run while keeping every
#    inherited property intact.
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("настройки для игры.", $true, $false, $false, $false, $false, $true, 1, $false, "настройки для игры. ", 2)

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastParaRange = $lastPara.Range
$trailingSpace = $d.Range($lastParaRange.End - 2, $lastParaRange.End - 1)
$trailingSpace.Font.Size = 20
$trailingSpace.Font.Size = 14

# ---------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark as a zero-length bookmark right at
#    the new end of the document (after the space we just appended).
#
#    Bookmarks.Add placed exactly at Content.End (or one character
#    before it) unexpectedly anchors at the start of the document, so a
#    throwaway placeholder character is appended first to push the
#    insertion point away from that edge case; the bookmark is added
#    just before the placeholder, which is then deleted again.
# ---------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertAfter("Z")

$afterPlaceholder = $d.Content.End
$bmPos = $afterPlaceholder - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tailEnd = $d.Content.End
$placeholderRange = $d.Range($tailEnd - 2, $tailEnd - 1)
$placeholderRange.Delete()
